$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$textRefs = @("D4", "D5", "D6", "D12", "D13", "D18", "D19", "D23", "D25", "D26", "D28", "D29", "D31", "D33", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D46", "D47", "D49")
foreach ($r in $textRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.610.40'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '2.563.74'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '578.10'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '143.26'
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '0.348'
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("D13").Value = '26.76'
$ws.Range("E13").Value = '  -4.21%  '
$ws.Range("D14").Value = '3.020.46'
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("D15").Value = '62.509.85'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = '2.551.21'
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '11.08'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("D19").Value = '337.58'
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '66.89'
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("E24").Value = '  -5.41%  '
$ws.Range("B25").Value = 'SuiNetwork'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D25").Value = '1.51'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  -4.64%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '7.92'
$ws.Range("E28").Value = '  -4.02%  '
$ws.Range("D29").Value = '8.17'
$ws.Range("E29").Value = '  -4.70%  '
$ws.Range("E30").Value = '  -2.72%  '
$ws.Range("D31").Value = '456.22'
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("D32").Value = '0.0₃0791'
$ws.Range("E32").Value = '  -4.81%  '
$ws.Range("D33").Value = '177.02'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = '0.394'
$ws.Range("E36").Value = '  -3.63%  '
$ws.Range("D37").Value = '18.78'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("D38").Value = '4.43'
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  -4.90%  '
$ws.Range("D41").Value = '40.49'
$ws.Range("E41").Value = '  +1.71%  '
$ws.Range("D42").Value = '157.07'
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("E43").Value = '  -4.90%  '
$ws.Range("D44").Value = '0.627'
$ws.Range("E44").Value = '  +2.87%  '
$ws.Range("D45").Value = '20.73'
$ws.Range("E45").Value = '  -3.75%  '
$ws.Range("D46").Value = '0.0531'
$ws.Range("E46").Value = '  -4.19%  '
$ws.Range("D47").Value = '0.0954'
$ws.Range("E47").Value = '  -2.55%  '
$ws.Range("E48").Value = '  -4.17%  '
$ws.Range("D49").Value = '17.87'
$ws.Range("E49").Value = '  -3.54%  '
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").Value = '  -5.29%  '
